$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Miguel Matias -> HEITOR MAMEDE (A2/B2 gain the underline style) ---
$ws.Range("A2").Value = "HEITOR MAMEDE"
$ws.Range("A2").Font.Underline = $true
$ws.Range("B2").Value = 83993238060
$ws.Range("B2").Font.Underline = $true
$ws.Range("C2").Value = "HEITOR MAMEDE.pdf"

# --- Row 3: Ramon -> ROBERTA MOURA (C3 loses the underline style) ---
$ws.Range("A3").Value = "ROBERTA MOURA"
$ws.Range("B3").Value = 83987317588
$ws.Range("C3").Value = "ROBERTA MOURA.pdf"
$ws.Range("C3").Font.Underline = $false

# --- Row 4: bella -> IAGO JULIANO (A4/B4 lose the underline style) ---
$ws.Range("A4").Value = "IAGO JULIANO"
$ws.Range("A4").Font.Underline = $false
$ws.Range("B4").Value = 83991071312
$ws.Range("B4").Font.Underline = $false
$ws.Range("C4").Value = "IAGO JULIANO.pdf"

# --- Row 5: new employee CAUA SANTANA ---
$ws.Range("A5").Value = "CAUA SANTANA"
$ws.Range("B5").Value = 83986901565
$ws.Range("B5").Font.Underline = $true
$ws.Range("C5").Value = "CAUA SANTANA.pdf"

# Row 6 (B6, underline-only placeholder cell) is left untouched.

# --- Stray formatted placeholder cells shift around ---
# B8 -> C8
$ws.Range("B8").Clear()
$ws.Range("C8").Font.Underline = $true

# C10 placeholder removed entirely
$ws.Range("C10").Clear()

# K23 -> K21
$ws.Range("K23").Clear()
$ws.Range("K21").Font.Underline = $true

# Matches the saved selection left behind in the source file (row 6 selected)
$ws.Rows(6).Select() | Out-Null
